# Trade #5 closed at 2026-02-17 07:52:39 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate metrics after the new trade closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.92    # Current Capital
$summary.Range("B4").Value = -0.08      # Total P&L $
$summary.Range("B5").Value = -0.32      # Total P&L %
$summary.Range("B6").Value = 5          # Total Trades
$summary.Range("B8").Value = 3          # Losing Trades
$summary.Range("B9").Value = 40         # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) reflects the new totals.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.92
$status.Range("D4").Value = 5
$status.Range("E4").Value = -0.08
$status.Range("F4").Value = -0.08
$status.Range("G4").Value = 40

# ---------------------------------------------------------------------------
# Helper: append the new trade row (row 6) to a trade-log sheet. Date/time
# text must stay text (not get auto-converted into date/time serials), so we
# force a text number format, assign the value, then drop back to the
# workbook's default "Normal" style so no stray formatting lingers.
# ---------------------------------------------------------------------------
function Add-TradeRow5($ws) {
    $ws.Range("A6").Value = 5

    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = "2026-02-17"
    $ws.Range("B6").Style = "Normal"

    $ws.Range("C6").NumberFormat = "@"
    $ws.Range("C6").Value = "07:52:32"
    $ws.Range("C6").Style = "Normal"

    $ws.Range("D6").Value = "MarketMaking"
    $ws.Range("E6").Value = "DOWN"
    $ws.Range("F6").Value = 0.3987
    $ws.Range("G6").Value = 0.383838
    $ws.Range("H6").Value = "CLOSED"
    $ws.Range("I6").Value = -3.7274
    $ws.Range("J6").Value = -0.01
    $ws.Range("K6").Value = 99.92
    $ws.Range("L6").Value = 0
    $ws.Range("M6").Value = 0
    $ws.Range("N6").Value = 0.6
    $ws.Range("O6").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P6").Value = "early_exit"
    $ws.Range("Q6").Value = 0.13
}

Add-TradeRow5 $wb.Worksheets.Item("All Trades")
Add-TradeRow5 $wb.Worksheets.Item("MarketMaking")
